$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 56, shifting existing rows 56-152 down to 57-153.
$ws.Rows.Item(56).Insert()

# Populate the newly inserted row 56 with the new record.
$ws.Cells.Item(56, 1).Value = 7
$ws.Cells.Item(56, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(56, 3).Value = "Ñuble"
$ws.Cells.Item(56, 4).Value = 45070
$ws.Cells.Item(56, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(56, 5).Value = 16
$ws.Cells.Item(56, 6).Value = "Fruta"
$ws.Cells.Item(56, 7).Value = 100103
$ws.Cells.Item(56, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(56, 9).Value = 100103002
$ws.Cells.Item(56, 10).Value = "Ciruela"
$ws.Cells.Item(56, 11).Value = "Angeleno"
$ws.Cells.Item(56, 12).Value = "Primera"
$ws.Cells.Item(56, 13).Value = 60
$ws.Cells.Item(56, 14).Value = 8000
$ws.Cells.Item(56, 15).Value = 9000
$ws.Cells.Item(56, 16).Value = 8500
$ws.Cells.Item(56, 17).Value = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(56, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(56, 19).Value = 472
$ws.Cells.Item(56, 20).Value = 18
